$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (target cluster changes from FAPs to ECs, plus value updates)
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.297418666666666
$ws.Range("H2").Value = 3.892256
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.162136
$ws.Range("N2").Value = 0.486408
$ws.Range("O2").Value = 0.0598694021675715
$ws.Range("P2").Value = 0.0598694021675715
$ws.Range("Q2").Value = 0.2103582729386666
$ws.Range("R2").Value = 1.893224456448
$ws.Range("S2").Value = 0.0598694021675715
$ws.Range("T2").Value = 0.0598694021675715

# Update existing row 3 (target cluster changes from M2 to FAPs, plus value updates)
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.297418666666666
$ws.Range("H3").Value = 3.892256
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.694965
$ws.Range("N3").Value = 5.084895
$ws.Range("O3").Value = 0.6258729785177741
$ws.Range("P3").Value = 0.6258729785177741
$ws.Range("Q3").Value = 2.199079230346666
$ws.Range("R3").Value = 19.79171307312
$ws.Range("S3").Value = 0.6258729785177741
$ws.Range("T3").Value = 0.6258729785177741

# New row 4 (target cluster M2)
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Il12a"
$ws.Range("C4").Value = "Il12rb1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.297418666666666
$ws.Range("H4").Value = 3.892256
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6978413333333333
$ws.Range("N4").Value = 2.093524
$ws.Range("O4").Value = 0.2576808570242738
$ws.Range("P4").Value = 0.2576808570242737
$ws.Range("Q4").Value = 0.9053923722382221
$ws.Range("R4").Value = 8.148531350143999
$ws.Range("S4").Value = 0.2576808570242738
$ws.Range("T4").Value = 0.2576808570242737

# New row 5 (target cluster sCs)
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Il12a"
$ws.Range("C5").Value = "Il12rb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.297418666666666
$ws.Range("H5").Value = 3.892256
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.153219
$ws.Range("N5").Value = 0.459657
$ws.Range("O5").Value = 0.05657676229038053
$ws.Range("P5").Value = 0.05657676229038053
$ws.Range("Q5").Value = 0.198789190688
$ws.Range("R5").Value = 1.789102716192
$ws.Range("S5").Value = 0.05657676229038053
$ws.Range("T5").Value = 0.05657676229038053
